# Apply the latest cryptocurrency price/volume snapshot to Sheet1.
# Columns: A=rank(unchanged) B=Coin C=Link D=Price(text) E=Volume(1h)(text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.031.35'
$ws.Range('E2').Value = '  +0.26%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.822.48'
$ws.Range('E3').Value = '  +0.26%  '
# Row 4
$ws.Range('E4').Value = '  +0.07%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.16'
$ws.Range('E5').Value = '  -0.14%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.619'
$ws.Range('E6').Value = '  +1.05%  '
# Row 7
$ws.Range('E7').Value = '  +0.10%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.13'
$ws.Range('E8').Value = '  -3.28%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.323'
$ws.Range('E9').Value = '  +5.40%  '
# Row 10
$ws.Range('E10').Value = '  +0.01%  '
# Row 11
$ws.Range('E11').Value = '  -0.80%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.084.33'
$ws.Range('E12').Value = '  +0.05%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.38'
$ws.Range('E13').Value = '  +3.05%  '
# Row 14
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.668'
$ws.Range('E14').Value = '  +1.63%  '
# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.807.07'
$ws.Range('E15').Value = '  -0.88%  '
# Row 16
$ws.Range('E16').Value = '  +0.06%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.050.91'
$ws.Range('E17').Value = '  +0.33%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.71'
$ws.Range('E18').Value = '  +0.49%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0787'
$ws.Range('E19').Value = '  +0.31%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.64'
$ws.Range('E20').Value = '  +1.23%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.08'
$ws.Range('E21').Value = '  +3.33%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.76'
$ws.Range('E22').Value = '  +2.70%  '
# Row 23
$ws.Range('E23').Value = '  +0.15%  '
# Row 24
$ws.Range('E24').Value = '  +1.14%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '174.05'
$ws.Range('E25').Value = '  +1.01%  '
# Row 26
$ws.Range('E26').Value = '  +0.84%  '
# Row 27
$ws.Range('E27').Value = '  +3.13%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.37'
$ws.Range('E28').Value = '  -0.01%  '
# Row 29
$ws.Range('E29').Value = '  -3.85%  '
# Row 30
$ws.Range('E30').Value = '  +0.07%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.01'
$ws.Range('E31').Value = '  +3.55%  '
# Row 32
$ws.Range('E32').Value = '  +0.73%  '
# Row 33
$ws.Range('E33').Value = '  +0.01%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.24'
$ws.Range('E34').Value = '  +11.74%  '
# Row 35
$ws.Range('E35').Value = '  +3.67%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.701'
$ws.Range('E36').Value = '  +4.01%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '93.22'
$ws.Range('E37').Value = '  +0.84%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.38'
$ws.Range('E38').Value = '  +8.02%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.342.34'
$ws.Range('E39').Value = '  +2.20%  '
# Row 40
$ws.Range('E40').Value = '  +1.12%  '
# Row 41
$ws.Range('E41').Value = '  +0.76%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.74'
$ws.Range('E42').Value = '  +0.67%  '
# Row 43
$ws.Range('E43').Value = '  -0.80%  '
# Row 44
$ws.Range('E44').Value = '  -0.92%  '
# Row 45
$ws.Range('E45').Value = '  -0.23%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.27'
$ws.Range('E46').Value = '  +1.26%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0521'
$ws.Range('E47').Value = '  +2.26%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.000.93'
$ws.Range('E48').Value = '  +0.24%  '
# Row 49
$ws.Range('E49').Value = '  +0.09%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0668'
$ws.Range('E50').Value = '  +4.65%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.21'
$ws.Range('E51').Value = '  +13.69%  '

Write-Output "Applied 84 cell updates"
